$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy style formatting from row 11 (A:H) down to row 12 first
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new attendance row 12: date 2020-03-06, time 9:40 AM
$ws.Range("A12").Value = 43896
$ws.Range("B12").Value = (9 + 40/60) / 24

# Update the active selection to reflect the new cursor position
$ws.Range("C13").Select()
